# Auto-generated script: update cached market-price / profit figures
# across all 8 job sheets, per scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 52632710
$ws.Range("I11").Value = 52632710
$ws.Range("K11").Value = 52632710
$ws.Range("M11").Value = -52632570
# Row 19
$ws.Range("H19").Value = 5913.05
$ws.Range("J19").Value = 983.5714
$ws.Range("L19").Value = 983.5714
$ws.Range("N19").Value = -1333.5714
# Row 43
$ws.Range("H43").Value = 550.6111
$ws.Range("I43").Value = 469.75
$ws.Range("J43").Value = 573.7143
$ws.Range("K43").Value = 469.75
$ws.Range("L43").Value = 573.7143
$ws.Range("M43").Value = -400.75
$ws.Range("N43").Value = -711.7143
# Row 62
$ws.Range("H62").Value = 3716.25
$ws.Range("I62").Value = 3399.5833
$ws.Range("J62").Value = 4191.25
$ws.Range("K62").Value = 3399.5833
$ws.Range("L62").Value = 4191.25
$ws.Range("M62").Value = -2775.5833
$ws.Range("N62").Value = -5439.25
# Row 64
$ws.Range("H64").Value = 3693.4805
$ws.Range("I64").Value = 3194.173
$ws.Range("J64").Value = 4732.04
$ws.Range("K64").Value = 3194.173
$ws.Range("L64").Value = 4732.04
$ws.Range("M64").Value = -2946.173
$ws.Range("N64").Value = -5228.04
# Row 65
$ws.Range("H65").Value = 3716.25
$ws.Range("I65").Value = 3399.5833
$ws.Range("J65").Value = 4191.25
$ws.Range("K65").Value = 16997.9165
$ws.Range("L65").Value = 20956.25
$ws.Range("M65").Value = -13877.9165
$ws.Range("N65").Value = -27196.25
# Row 67
$ws.Range("H67").Value = 3693.4805
$ws.Range("I67").Value = 3194.173
$ws.Range("J67").Value = 4732.04
$ws.Range("K67").Value = 3194.173
$ws.Range("L67").Value = 4732.04
$ws.Range("M67").Value = -2336.173
$ws.Range("N67").Value = -6448.04
# Row 98
$ws.Range("H98").Value = 863.3889
$ws.Range("I98").Value = 929.6667
$ws.Range("J98").Value = 532
$ws.Range("K98").Value = 929.6667
$ws.Range("L98").Value = 532
$ws.Range("M98").Value = 568.3333
$ws.Range("N98").Value = -3528
# Row 112
$ws.Range("H112").Value = 1840
$ws.Range("J112").Value = 6000
$ws.Range("L112").Value = 18000
$ws.Range("N112").Value = -20216
# Row 122
$ws.Range("H122").Value = 863.3889
$ws.Range("I122").Value = 929.6667
$ws.Range("J122").Value = 532
$ws.Range("K122").Value = 2789.0001
$ws.Range("L122").Value = 1596
$ws.Range("M122").Value = -339.0001000000002
$ws.Range("N122").Value = -6496
# Row 132
$ws.Range("H132").Value = 723.63635
$ws.Range("I132").Value = 517.55554
$ws.Range("J132").Value = 1651
$ws.Range("K132").Value = 1552.66662
$ws.Range("L132").Value = 4953
$ws.Range("M132").Value = 977.33338
$ws.Range("N132").Value = -10013

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2503.9707
$ws.Range("I61").Value = 2274.8696
$ws.Range("J61").Value = 2983
$ws.Range("K61").Value = 2274.8696
$ws.Range("L61").Value = 2983
$ws.Range("M61").Value = -2062.8696
$ws.Range("N61").Value = -3407
# Row 136
$ws.Range("H136").Value = 2503.9707
$ws.Range("I136").Value = 2274.8696
$ws.Range("J136").Value = 2983
$ws.Range("K136").Value = 6824.6088
$ws.Range("L136").Value = 8949
$ws.Range("M136").Value = -4274.6088
$ws.Range("N136").Value = -14049

$ws = $wb.Worksheets.Item("BSM")
# Row 17
$ws.Range("H17").Value = 495
$ws.Range("I17").Value = 495
$ws.Range("K17").Value = 495
$ws.Range("M17").Value = -323
# Row 21
$ws.Range("H21").Value = 13736.2
$ws.Range("J21").Value = 13736.2
$ws.Range("L21").Value = 13736.2
$ws.Range("N21").Value = -14208.2
# Row 92
$ws.Range("H92").Value = 1500
$ws.Range("J92").Value = 1500
$ws.Range("L92").Value = 1500
$ws.Range("N92").Value = -6492

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 864.0714
$ws.Range("I7").Value = 1513.1428
$ws.Range("K7").Value = 1513.1428
$ws.Range("M7").Value = -1400.1428
# Row 31
$ws.Range("H31").Value = 1335.27
$ws.Range("I31").Value = 644.5577
$ws.Range("J31").Value = 2083.5417
$ws.Range("K31").Value = 644.5577
$ws.Range("L31").Value = 2083.5417
$ws.Range("M31").Value = -349.5577
$ws.Range("N31").Value = -2673.5417
# Row 34
$ws.Range("H34").Value = 1335.27
$ws.Range("I34").Value = 644.5577
$ws.Range("J34").Value = 2083.5417
$ws.Range("K34").Value = 644.5577
$ws.Range("L34").Value = 2083.5417
$ws.Range("M34").Value = -442.5577
$ws.Range("N34").Value = -2487.5417
# Row 62
$ws.Range("H62").Value = 4275.5
$ws.Range("J62").Value = 4327
$ws.Range("L62").Value = 4327
$ws.Range("N62").Value = -5575
# Row 65
$ws.Range("H65").Value = 4275.5
$ws.Range("J65").Value = 4327
$ws.Range("L65").Value = 21635
$ws.Range("N65").Value = -27875
# Row 134
$ws.Range("H134").Value = 3198.9556
$ws.Range("I134").Value = 4512.84
$ws.Range("J134").Value = 1556.6
$ws.Range("K134").Value = 13538.52
$ws.Range("L134").Value = 4669.799999999999
$ws.Range("M134").Value = -11003.52
$ws.Range("N134").Value = -9739.799999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 308.2857
$ws.Range("I2").Value = 422.6
$ws.Range("J2").Value = 22.5
$ws.Range("K2").Value = 2535.6
$ws.Range("L2").Value = 135
$ws.Range("M2").Value = -2422.6
$ws.Range("N2").Value = -361
# Row 12
$ws.Range("H12").Value = 60.95
$ws.Range("J12").Value = 85.07143000000001
$ws.Range("L12").Value = 255.21429
$ws.Range("N12").Value = -601.21429
# Row 28
$ws.Range("H28").Value = 397.5
$ws.Range("I28").Value = 397.5
$ws.Range("K28").Value = 1192.5
$ws.Range("M28").Value = -960.5
# Row 33
$ws.Range("H33").Value = 119.68
$ws.Range("I33").Value = 68.82353000000001
$ws.Range("J33").Value = 227.75
$ws.Range("K33").Value = 412.94118
$ws.Range("L33").Value = 1366.5
$ws.Range("M33").Value = -129.94118
$ws.Range("N33").Value = -1932.5
# Row 45
$ws.Range("H45").Value = 593.2
$ws.Range("J45").Value = 641.5
$ws.Range("L45").Value = 1924.5
$ws.Range("N45").Value = -2988.5
# Row 46
$ws.Range("H46").Value = 1758.1538
$ws.Range("J46").Value = 2250
$ws.Range("L46").Value = 6750
$ws.Range("N46").Value = -6932
# Row 55
$ws.Range("H55").Value = 8085.615
$ws.Range("J55").Value = 8085.615
$ws.Range("L55").Value = 24256.845
$ws.Range("N55").Value = -24610.845
# Row 64
$ws.Range("H64").Value = 1300.5264
$ws.Range("I64").Value = 500
$ws.Range("J64").Value = 1345
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 4035
$ws.Range("M64").Value = -1230
$ws.Range("N64").Value = -4575
# Row 67
$ws.Range("H67").Value = 1300.5264
$ws.Range("I67").Value = 500
$ws.Range("J67").Value = 1345
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 4035
$ws.Range("M67").Value = -564
$ws.Range("N67").Value = -5907
# Row 87
$ws.Range("H87").Value = 18400.5
$ws.Range("I87").Value = 2077
$ws.Range("J87").Value = 34724
$ws.Range("K87").Value = 6231
$ws.Range("L87").Value = 104172
$ws.Range("M87").Value = -4983
$ws.Range("N87").Value = -106668
# Row 90
$ws.Range("H90").Value = 18400.5
$ws.Range("I90").Value = 2077
$ws.Range("J90").Value = 34724
$ws.Range("K90").Value = 18693
$ws.Range("L90").Value = 312516
$ws.Range("M90").Value = -12453
$ws.Range("N90").Value = -324996
# Row 107
$ws.Range("H107").Value = 800.2954999999999
$ws.Range("I107").Value = 885.26086
$ws.Range("J107").Value = 707.2381
$ws.Range("K107").Value = 2655.78258
$ws.Range("L107").Value = 2121.7143
$ws.Range("M107").Value = -735.7825800000001
$ws.Range("N107").Value = -5961.7143
# Row 109
$ws.Range("H109").Value = 3694.1177
$ws.Range("I109").Value = 1081.6
$ws.Range("J109").Value = 4782.6665
$ws.Range("K109").Value = 3244.8
$ws.Range("L109").Value = 14347.9995
$ws.Range("M109").Value = -2204.8
$ws.Range("N109").Value = -16427.9995

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3268.1365
$ws.Range("I80").Value = 2992.3076
$ws.Range("J80").Value = 3666.5557
$ws.Range("K80").Value = 2992.3076
$ws.Range("L80").Value = 3666.5557
$ws.Range("M80").Value = -1994.3076
$ws.Range("N80").Value = -5662.5557
# Row 83
$ws.Range("H83").Value = 3268.1365
$ws.Range("I83").Value = 2992.3076
$ws.Range("J83").Value = 3666.5557
$ws.Range("K83").Value = 14961.538
$ws.Range("L83").Value = 18332.7785
$ws.Range("M83").Value = -9969.538
$ws.Range("N83").Value = -28316.7785

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 55200.844
$ws.Range("I7").Value = 65154.125
$ws.Range("J7").Value = 2116.6667
$ws.Range("K7").Value = 65154.125
$ws.Range("L7").Value = 2116.6667
$ws.Range("M7").Value = -65042.125
$ws.Range("N7").Value = -2340.6667
# Row 126
$ws.Range("H126").Value = 55200.844
$ws.Range("I126").Value = 65154.125
$ws.Range("J126").Value = 2116.6667
$ws.Range("K126").Value = 195462.375
$ws.Range("L126").Value = 6350.000100000001
$ws.Range("M126").Value = -192992.375
$ws.Range("N126").Value = -11290.0001
# Row 132
$ws.Range("H132").Value = 5001
$ws.Range("I132").Value = 4500.9165
$ws.Range("J132").Value = 7001.3335
$ws.Range("K132").Value = 13502.7495
$ws.Range("L132").Value = 21004.0005
$ws.Range("M132").Value = -10972.7495
$ws.Range("N132").Value = -26064.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 122
$ws.Range("H122").Value = 34153.84
$ws.Range("I122").Value = 60579.707
$ws.Range("J122").Value = 2065.2856
$ws.Range("K122").Value = 181739.121
$ws.Range("L122").Value = 6195.8568
$ws.Range("M122").Value = -179289.121
$ws.Range("N122").Value = -11095.8568
